$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price and Volume columns to be treated as text so values such as
# "1.007" or "28.274.68" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '28.274.68'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.885.58'
$ws.Range("E3").Value = '  +0.83%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '314.13'
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.5141'
$ws.Range("E7").Value = '  +1.19%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3911'
$ws.Range("E8").Value = '  +2.67%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.08384'
$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.120'
$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").Value = '6.247'
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '20.78'
$ws.Range("E12").Value = '  +1.48%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.875.16'
$ws.Range("E13").Value = '  +0.25%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '7.312'
$ws.Range("E14").Value = '  +1.73%  '

$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = '1.006'
$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.00001109'
$ws.Range("E16").Value = '  +1.12%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '91.58'
$ws.Range("E17").Value = '  +0.95%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.06673'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '17.85'
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.068'
$ws.Range("E21").Value = '  +0.95%  '

$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '28.306.58'
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '11.21'
$ws.Range("E23").Value = '  +0.88%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '2.266'
$ws.Range("E24").Value = '  +0.14%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.094.01'
$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.527'
$ws.Range("E26").Value = '  -1.99%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '158.98'
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '20.69'
$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '125.55'
$ws.Range("E29").Value = '  -0.33%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.1066'
$ws.Range("E30").Value = '  +0.92%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.046'
$ws.Range("E31").Value = '  +0.30%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.888'
$ws.Range("E32").Value = '  +5.21%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '3.597'
$ws.Range("E33").Value = '  +0.21%  '

$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").Value = '9.820'
$ws.Range("E34").Value = '  +0.87%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '0.02460'
$ws.Range("E35").Value = '  +0.69%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.06579'
$ws.Range("E36").Value = '  +0.25%  '

$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2197'
$ws.Range("E37").Value = '  +1.52%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '1.214'
$ws.Range("E38").Value = '  +0.41%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.6539'
$ws.Range("E39").Value = '  +1.45%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '5.042'
$ws.Range("E40").Value = '  +3.35%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.232'
$ws.Range("E41").Value = '  -0.37%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '11.32'
$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").Value = '0.6134'
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '13.20'
$ws.Range("E44").Value = '  +1.77%  '

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '1.295'
$ws.Range("E45").Value = '  +0.48%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.684'
$ws.Range("E46").Value = '  +0.60%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.027'
$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = '1.236'
$ws.Range("E48").Value = '  +1.85%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '121.91'
$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '78.90'
$ws.Range("E50").Value = '  -1.56%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.06918'
$ws.Range("E51").Value = '  +0.80%  '
